# Update the movie-information workbook:
#   - A2 (Director's name) changes from "Steven Spielberg" to "ur mom"
#   - Active selection on Sheet1 moves from A2 to A3

try {
    $wb = $excel.ActiveWorkbook
    $ws = $wb.ActiveSheet

    $ws.Range("A2").Value = "ur mom"
    $ws.Range("A3").Select()
}
catch {
    Write-Output "Error while editing workbook: $_"
}
